$d = $word.ActiveDocument

# Run properties shared by every run in these two "[<MSSVn - Ho va ten sinh
# vien>]" placeholder paragraphs (Arial, blue, 15pt).
$rPr = '<w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="0000FF"/><w:sz w:val="30"/><w:szCs w:val="30"/><w:lang w:val="en-US"/></w:rPr>'

function New-PkgXml($innerRunsXml) {
    return '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>' + $innerRunsXml + '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
}

# Resolve the two target paragraphs (by their current placeholder text) up
# front, before mutating anything, so the replacement loop below never walks
# a collection that is being edited underneath it.
$mssv2Para = $null
$mssv4Para = $null
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -match "^\[<MSSV2") { $mssv2Para = $p }
    elseif ($t -match "^\[<MSSV4") { $mssv4Para = $p }
}

if ($mssv2Para -ne $null) {
    $rng = $mssv2Para.Range
    $target = $d.Range($rng.Start, $rng.End - 1)
    $runs = '<w:r>' + $rPr + '<w:t>22850213</w:t></w:r>' +
            '<w:r>' + $rPr + '<w:t xml:space="preserve">&#8211; </w:t></w:r>' +
            '<w:r>' + $rPr + '<w:t>Phan Thi&#234;n Qu&#7889;c</w:t></w:r>'
    $target.InsertXML((New-PkgXml $runs))
}

if ($mssv4Para -ne $null) {
    $rng = $mssv4Para.Range
    $target = $d.Range($rng.Start, $rng.End - 1)
    $runs = '<w:r>' + $rPr + '<w:t>22810209</w:t></w:r>' +
            '<w:r>' + $rPr + '<w:t xml:space="preserve">&#8211; </w:t></w:r>' +
            '<w:r>' + $rPr + '<w:t>L&#226;m Tr&#7885;ng Ngh&#297;a</w:t></w:r>'
    $target.InsertXML((New-PkgXml $runs))
}
